# Auto-generated script to apply scheduled-runner price refresh to Ifrit_Profits workbook
# Updates cached (static) market-price/profit figures across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 675.3333
$ws.Range("I99").Value = 468.2857
$ws.Range("K99").Value = 1404.8571
$ws.Range("M99").Value = 93.14289999999983
$ws.Range("H100").Value = 1459.5555
$ws.Range("I100").Value = 1367.1111
$ws.Range("J100").Value = 1644.4445
$ws.Range("K100").Value = 1367.1111
$ws.Range("L100").Value = 1644.4445
$ws.Range("M100").Value = -826.1111000000001
$ws.Range("N100").Value = -2726.4445
$ws.Range("H113").Value = 2120
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = 254
$ws.Range("H137").Value = 23258546
$ws.Range("I137").Value = 1663.1177
$ws.Range("J137").Value = 38464970
$ws.Range("K137").Value = 4989.3531
$ws.Range("L137").Value = 115394910
$ws.Range("M137").Value = -2439.3531
$ws.Range("N137").Value = -115400010
$ws.Range("H138").Value = 3157.0566
$ws.Range("I138").Value = 2698.3242
$ws.Range("J138").Value = 4217.875
$ws.Range("K138").Value = 8094.9726
$ws.Range("L138").Value = 12653.625
$ws.Range("M138").Value = -2954.9726
$ws.Range("N138").Value = -22933.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7356936
$ws.Range("I2").Value = 6490
$ws.Range("J2").Value = 14707382
$ws.Range("K2").Value = 6490
$ws.Range("L2").Value = 14707382
$ws.Range("M2").Value = -6377
$ws.Range("N2").Value = -14707608
$ws.Range("H37").Value = 8558.200000000001
$ws.Range("I37").Value = 1484
$ws.Range("J37").Value = 10326.75
$ws.Range("K37").Value = 1484
$ws.Range("L37").Value = 10326.75
$ws.Range("M37").Value = -1211
$ws.Range("N37").Value = -10872.75
$ws.Range("H61").Value = 3475759.5
$ws.Range("I61").Value = 5294029
$ws.Range("J61").Value = 4517.636
$ws.Range("K61").Value = 5294029
$ws.Range("L61").Value = 4517.636
$ws.Range("M61").Value = -5293817
$ws.Range("N61").Value = -4941.636
$ws.Range("H102").Value = 2672.4
$ws.Range("I102").Value = 2812.5557
$ws.Range("J102").Value = 1411
$ws.Range("K102").Value = 2812.5557
$ws.Range("L102").Value = 1411
$ws.Range("M102").Value = -1190.5557
$ws.Range("N102").Value = -4655
$ws.Range("H116").Value = 7356936
$ws.Range("I116").Value = 6490
$ws.Range("J116").Value = 14707382
$ws.Range("K116").Value = 6490
$ws.Range("L116").Value = 14707382
$ws.Range("M116").Value = -4196
$ws.Range("N116").Value = -14711970
$ws.Range("H132").Value = 1233905.6
$ws.Range("I132").Value = 1726106.9
$ws.Range("J132").Value = 3402.5
$ws.Range("K132").Value = 5178320.699999999
$ws.Range("L132").Value = 10207.5
$ws.Range("M132").Value = -5175790.699999999
$ws.Range("N132").Value = -15267.5
$ws.Range("H136").Value = 3475759.5
$ws.Range("I136").Value = 5294029
$ws.Range("J136").Value = 4517.636
$ws.Range("K136").Value = 15882087
$ws.Range("L136").Value = 13552.908
$ws.Range("M136").Value = -15879537
$ws.Range("N136").Value = -18652.908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7356936
$ws.Range("I3").Value = 6490
$ws.Range("J3").Value = 14707382
$ws.Range("K3").Value = 6490
$ws.Range("L3").Value = 14707382
$ws.Range("M3").Value = -6376
$ws.Range("N3").Value = -14707610
$ws.Range("H94").Value = 258.33334
$ws.Range("I94").Value = 258.33334
$ws.Range("K94").Value = 258.33334
$ws.Range("M94").Value = 192.66666
$ws.Range("H105").Value = 1622.3077
$ws.Range("I105").Value = 1409
$ws.Range("J105").Value = 2333.3333
$ws.Range("K105").Value = 1409
$ws.Range("L105").Value = 2333.3333
$ws.Range("M105").Value = 338
$ws.Range("N105").Value = -5827.3333
$ws.Range("H134").Value = 28588644
$ws.Range("I134").Value = 33352998
$ws.Range("J134").Value = 2514
$ws.Range("K134").Value = 100058994
$ws.Range("L134").Value = 7542
$ws.Range("M134").Value = -100056459
$ws.Range("N134").Value = -12612

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2960.3076
$ws.Range("I31").Value = 1401.2632
$ws.Range("J31").Value = 7192
$ws.Range("K31").Value = 1401.2632
$ws.Range("L31").Value = 7192
$ws.Range("M31").Value = -1106.2632
$ws.Range("N31").Value = -7782
$ws.Range("H34").Value = 2960.3076
$ws.Range("I34").Value = 1401.2632
$ws.Range("J34").Value = 7192
$ws.Range("K34").Value = 1401.2632
$ws.Range("L34").Value = 7192
$ws.Range("M34").Value = -1199.2632
$ws.Range("N34").Value = -7596
$ws.Range("H69").Value = 10551.714
$ws.Range("I69").Value = 7310.3335
$ws.Range("K69").Value = 7310.3335
$ws.Range("M69").Value = -6561.3335
$ws.Range("H72").Value = 10551.714
$ws.Range("I72").Value = 7310.3335
$ws.Range("K72").Value = 21931.0005
$ws.Range("M72").Value = -18187.0005
$ws.Range("H107").Value = 2315909.2
$ws.Range("I107").Value = 3206145.2
$ws.Range("J107").Value = 1296
$ws.Range("K107").Value = 3206145.2
$ws.Range("L107").Value = 1296
$ws.Range("M107").Value = -3204225.2
$ws.Range("N107").Value = -5136
$ws.Range("H132").Value = 2111.8484
$ws.Range("I132").Value = 1880.5172
$ws.Range("K132").Value = 5641.5516
$ws.Range("M132").Value = -3111.5516

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2087.5
$ws.Range("I97").Value = 1928.5714
$ws.Range("J97").Value = 3200
$ws.Range("K97").Value = 1928.5714
$ws.Range("L97").Value = 3200
$ws.Range("M97").Value = -1432.5714
$ws.Range("N97").Value = -4192
$ws.Range("H132").Value = 2310.7097
$ws.Range("I132").Value = 2019.037
$ws.Range("J132").Value = 4279.5
$ws.Range("K132").Value = 6057.111
$ws.Range("L132").Value = 12838.5
$ws.Range("M132").Value = -3527.111
$ws.Range("N132").Value = -17898.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2250
$ws.Range("J61").Value = 2400
$ws.Range("L61").Value = 2400
$ws.Range("N61").Value = -2804
$ws.Range("H113").Value = 2250
$ws.Range("J113").Value = 2400
$ws.Range("L113").Value = 2400
$ws.Range("N113").Value = -6740

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7835.65
$ws.Range("I136").Value = 10170.214
$ws.Range("J136").Value = 2388.3333
$ws.Range("K136").Value = 30510.642
$ws.Range("L136").Value = 7164.999899999999
$ws.Range("M136").Value = -27960.642
$ws.Range("N136").Value = -12264.9999
